$p = $ppt.ActivePresentation

# --- Slide 2, Content Placeholder, paragraph 3:
# "...which may cause the some task durations to blow out."
#   -> "...which may cause some task durations to blow out."
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(3, 1)
$run2 = $para2.Runs(1, 1)
$run2.Text = "Very little is known of the target company’s infrastructure other than a basic device count and the fact that a company website exists. This has meant that I have, if anything, deliberately overestimated the duration of some tasks. There may yet be a few surprises in store which may cause some task durations to blow out."

# --- Slide 3, Content Placeholder, paragraph 4:
# Collapse the three runs ("I need to monitor actual vs " / "planned task " /
# "durations on a regular basis and update the project plan as required.")
# into a single run of text.
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange

$para3c = $tr3.Paragraphs(4, 1)
$run3c = $para3c.Runs(3, 1)
$run3c.Text = ""

$para3b = $tr3.Paragraphs(4, 1)
$run3b = $para3b.Runs(2, 1)
$run3b.Text = ""

$para3a = $tr3.Paragraphs(4, 1)
$run3a = $para3a.Runs(1, 1)
$run3a.Text = "I need to monitor actual vs planned task durations on a regular basis and update the project plan as required."
